$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 new match data (24/10/2025 Al-Ittihad vs Al-Hilal)
$ws.Range("A7").Value = "24/10/2025"
$ws.Range("B7").Value = "Al-Ittihad"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Al-Hilal"
$ws.Range("F7").Value = "W"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1.72
$ws.Range("L7").Value = 0.83
$ws.Range("M7").Value = 11
$ws.Range("N7").Value = 15
$ws.Range("O7").Value = 5
$ws.Range("P7").Value = 4
